# Update "想去人数" (want-to-go count) values in column F on the
# "展览", "演出", and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$changes = @{
    "展览"     = @{ 2 = 66; 4 = 2038; 5 = 332; 7 = 95; 9 = 10506; 14 = 407; 15 = 7401; 16 = 1111; 18 = 196; 20 = 3300 }
    "演出"     = @{ 2 = 20 }
    "全部类型" = @{ 2 = 66; 4 = 2038; 5 = 332; 7 = 20; 8 = 95; 12 = 10506; 17 = 407; 18 = 7401; 19 = 1111; 21 = 196; 23 = 3300 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $changes[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}
